# Generate Report for Handback
# - Update Status text ("Ready for handoff" -> "Handed back: in sync with en-US")
# - Set the Latest Handback DateTime for zh-cn / de-de
# - Fill in "Latest Target File" (F) / "Latest Handback File" (G) columns with
#   hyperlinked file names on both the zh-cn and de-de sheets

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Update the status text everywhere it appears (Overview + zh-cn + de-de)
# ---------------------------------------------------------------------------
$statusCells = @{
    "Overview" = @("B2", "C2", "B3", "C3")
    "zh-cn"    = @("C2", "C3")
    "de-de"    = @("C2", "C3")
}

foreach ($sheetName in $statusCells.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $statusCells[$sheetName]) {
        $cell = $ws.Range($addr)
        if ($cell.Text -eq $statusOld) {
            $cell.Value = $statusNew
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Update "Latest Handback DateTime" (column H) on zh-cn and de-de
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-03-24 10:18:40"
$wsZhCn.Range("H3").Value = "2016-03-24 10:18:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-03-24 10:18:47"
$wsDeDe.Range("H3").Value = "2016-03-24 10:18:47"

# ---------------------------------------------------------------------------
# 3. Fill in F (Latest Target File) and G (Latest Handback File) columns,
#    with hyperlinks matching the existing Source File Name (A) / Latest
#    Handoff File (D) links for row 2 on both sheets.
# ---------------------------------------------------------------------------
function Add-ReportLinks($ws, $mdFileName, $xlfFileName) {
    # Grab the URL already used for the row-2 hyperlinks so the new ones
    # point at the same targets.
    $mdUrl = ""
    $xlfUrl = ""
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq "`$A`$2") { $mdUrl = $hl.Address }
        if ($addr -eq "`$D`$2") { $xlfUrl = $hl.Address }
    }

    $ws.Range("F2").Value = $mdFileName
    $ws.Range("G2").Value = $xlfFileName
    $ws.Range("F3").Value = $mdFileName
    $ws.Range("G3").Value = $xlfFileName

    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName)
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $xlfFileName)
    $ws.Hyperlinks.Add($ws.Range("F3"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName)
    $ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $xlfFileName)
}

Add-ReportLinks $wsZhCn "437cfbd5-767f-4178-a01b-f91116985aef.md" "437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.zh-cn.xlf"

Add-ReportLinks $wsDeDe "437cfbd5-767f-4178-a01b-f91116985aef.md" "437cfbd5-767f-4178-a01b-f91116985aef.d8c2c853bf6105a22b55831bbb4018e9b23b0fbb.de-de.xlf"

Write-Host "Report for handback generated."
